$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    if ($value -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") "67.384.87"
$ws.Range("E2").Value = "  -0.77%  "

Set-TextValue $ws.Range("D3") "3.522.50"
$ws.Range("E3").Value = "  -2.94%  "

$ws.Range("E4").Value = "  +0.32%  "

Set-TextValue $ws.Range("D5") "202.80"
$ws.Range("E5").Value = "  +3.82%  "

Set-TextValue $ws.Range("D6") "552.73"
$ws.Range("E6").Value = "  -5.76%  "

Set-TextValue $ws.Range("D7") "3.514.92"
$ws.Range("E7").Value = "  -3.03%  "

$ws.Range("E8").Value = "  -2.43%  "

Set-TextValue $ws.Range("D10") "0.662"
$ws.Range("E10").Value = "  -3.37%  "

Set-TextValue $ws.Range("D11") "61.81"
$ws.Range("E11").Value = "  +11.29%  "

$ws.Range("E12").Value = "  -5.47%  "

Set-TextValue $ws.Range("D13") "0.0000276"
$ws.Range("E13").Value = "  -2.57%  "

Set-TextValue $ws.Range("D14") "9.89"
$ws.Range("E14").Value = "  -1.65%  "

Set-TextValue $ws.Range("D15") "4.087.48"
$ws.Range("E15").Value = "  -2.35%  "

Set-TextValue $ws.Range("D16") "3.512.99"
$ws.Range("E16").Value = "  -2.88%  "

$ws.Range("E17").Value = "  -1.06%  "

$ws.Range("E18").Value = "  +0.45%  "

Set-TextValue $ws.Range("D19") "67.137.78"
$ws.Range("E19").Value = "  -0.94%  "

Set-TextValue $ws.Range("D20") "11.94"
$ws.Range("E20").Value = "  -4.42%  "

$ws.Range("E21").Value = "  -4.14%  "

Set-TextValue $ws.Range("D22") "391.12"
$ws.Range("E22").Value = "  -3.62%  "

$ws.Range("E23").Value = "  -5.27%  "

$ws.Range("E24").Value = "  -9.55%  "

Set-TextValue $ws.Range("D25") "82.78"
$ws.Range("E25").Value = "  -4.12%  "

$ws.Range("E26").Value = "  -4.62%  "

$ws.Range("E27").Value = "  -4.97%  "

Set-TextValue $ws.Range("D28") "3.74"
$ws.Range("E28").Value = "  -3.77%  "

Set-TextValue $ws.Range("D29") "8.94"
$ws.Range("E29").Value = "  -3.56%  "

Set-TextValue $ws.Range("D30") "30.85"
$ws.Range("E30").Value = "  -2.31%  "

Set-TextValue $ws.Range("D31") "7.38"
$ws.Range("E31").Value = "  -9.88%  "

Set-TextValue $ws.Range("D32") "681.23"
$ws.Range("E32").Value = "  +1.62%  "

Set-TextValue $ws.Range("D33") "11.80"
$ws.Range("E33").Value = "  -4.10%  "

Set-TextValue $ws.Range("D34") "63.64"
$ws.Range("E34").Value = "  -1.45%  "

Set-TextValue $ws.Range("D35") "0.112"
$ws.Range("E35").Value = "  -5.48%  "

Set-TextValue $ws.Range("D36") "39.89"
$ws.Range("E36").Value = "  -7.75%  "

$ws.Range("E37").Value = "  -2.82%  "

Set-TextValue $ws.Range("D38") "0.999"

Set-TextValue $ws.Range("D39") "3.08"
$ws.Range("E39").Value = "  -0.91%  "

Set-TextValue $ws.Range("D40") "3.115.52"
$ws.Range("E40").Value = "  -1.99%  "

Set-TextValue $ws.Range("D41") "0.131"
$ws.Range("E41").Value = "  -2.32%  "

Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  +0.07%  "

Set-TextValue $ws.Range("D43") "0.0₃0707"
$ws.Range("E43").Value = "  -11.99%  "

$ws.Range("E46").Value = "  +6.23%  "

Set-TextValue $ws.Range("D47") "0.0402"
$ws.Range("E47").Value = "  -4.71%  "

Set-TextValue $ws.Range("D48") "2.97"
$ws.Range("E48").Value = "  -4.86%  "

$ws.Range("E49").Value = "  -3.45%  "

Set-TextValue $ws.Range("D50") "137.16"
$ws.Range("E50").Value = "  -4.32%  "

Set-TextValue $ws.Range("D51") "8.33"
$ws.Range("E51").Value = "  -6.41%  "

# Row 44 and 45 content swap (Fetch.AI / dogwifhat reorder)
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D44") "2.58"
$ws.Range("E44").Value = "  -12.72%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D45") "2.84"
$ws.Range("E45").Value = "  +12.53%  "